$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-7 (old second gene triple, MuSCs as sender) - remaining rows shift up
$ws.Rows("5:7").Delete()

# Update A2:D4 (sending cluster / ligand / receptor / target cluster) for the
# new TPM-derived pairing: sender is now ECs/MuSCs/C1qb, target cycles
# through Lrp1, FAPs, ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "MuSCs"
$ws.Range("C2").Value = "C1qb"
$ws.Range("D2").Value = "Lrp1"

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "MuSCs"
$ws.Range("C3").Value = "C1qb"
$ws.Range("D3").Value = "FAPs"

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "MuSCs"
$ws.Range("C4").Value = "C1qb"
$ws.Range("D4").Value = "ECs"

# Row 2 numeric columns E..T
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.100566
$ws.Range("H2").Value = 0.301698
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 0.347582779512
$ws.Range("R2").Value = 3.128245015608
$ws.Range("S2").Value = 0.009841535807677501
$ws.Range("T2").Value = 0.0098415358076775

# Row 3 numeric columns E..T
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.100566
$ws.Range("H3").Value = 0.301698
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 30.330720550812
$ws.Range("R3").Value = 272.976484957308
$ws.Range("S3").Value = 0.8587907398420774
$ws.Range("T3").Value = 0.8587907398420773

# Row 4 numeric columns E..T
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.100566
$ws.Range("H4").Value = 0.301698
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 4.639637517978001
$ws.Range("R4").Value = 41.75673766180201
$ws.Range("S4").Value = 0.1313677243502452
$ws.Range("T4").Value = 0.1313677243502452
